$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new field was captured for this document type: "Bimestre o Período Anual".
# Append it as a new data row at the bottom of the existing table (row 35),
# following the same 5-column layout as the rest of the rows:
#   A: Concepto   B: Campo   C: Valor   D: Año   E: Razón Social
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Copy formatting from the last existing data row so the new row matches style.
$ws.Range("A" + ($newRow - 1) + ":E" + ($newRow - 1)).Copy() | Out-Null
$ws.Range("A" + $newRow + ":E" + $newRow).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = "Bimestre o Período Anual"
$ws.Cells.Item($newRow, 2).Value = "N/A"
$ws.Cells.Item($newRow, 3).Value = "Anual"
$ws.Cells.Item($newRow, 4).Value = 2019
$ws.Cells.Item($newRow, 5).Value = "INVERSIONES ORTIZ VASQUEZ HERMANOS S A S"
